$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift data rows 98-182 down by one record (new row inserted at the top of the shifted block) ---
$ws.Range("D98").Value = 44587
$ws.Range("D99").Value = 44586
$ws.Range("J99").Value = 900
$ws.Range("K99").Value = 300
$ws.Range("M99").Value = 300
$ws.Range("P99").Value = 300
$ws.Range("D100").Value = 44252
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 280
$ws.Range("L100").Value = 300
$ws.Range("M100").Value = 290
$ws.Range("O100").Value = 'Región del Maule'
$ws.Range("P100").Value = 290
$ws.Range("D101").Value = 44550
$ws.Range("J101").Value = 900
$ws.Range("K101").Value = 600
$ws.Range("L101").Value = 600
$ws.Range("M101").Value = 600
$ws.Range("O101").Value = 'Región de O''Higgins'
$ws.Range("P101").Value = 600
$ws.Range("D102").Value = 44524
$ws.Range("H102").Value = 'Camote'
$ws.Range("I102").Value = '1a nueva(o)'
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 700
$ws.Range("L102").Value = 700
$ws.Range("M102").Value = 700
$ws.Range("O102").Value = 'Perú'
$ws.Range("P102").Value = 700
$ws.Range("D103").Value = 44487
$ws.Range("H103").Value = 'Paine'
$ws.Range("I103").Value = '1a (guarda)'
$ws.Range("J103").Value = 3000
$ws.Range("K103").Value = 100
$ws.Range("L103").Value = 100
$ws.Range("M103").Value = 100
$ws.Range("P103").Value = 100
$ws.Range("D104").Value = 44267
$ws.Range("I104").Value = '1a (cosecha)'
$ws.Range("J104").Value = 800
$ws.Range("D105").Value = 44327
$ws.Range("I105").Value = '1a (guarda)'
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 300
$ws.Range("M105").Value = 300
$ws.Range("P105").Value = 300
$ws.Range("D106").Value = 44582
$ws.Range("J106").Value = 700
$ws.Range("K106").Value = 250
$ws.Range("M106").Value = 271
$ws.Range("P106").Value = 271
$ws.Range("D107").Value = 44235
$ws.Range("J107").Value = 800
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 300
$ws.Range("P107").Value = 300
$ws.Range("D108").Value = 44580
$ws.Range("H108").Value = 'Camote'
$ws.Range("I108").Value = '1a nueva(o)'
$ws.Range("J108").Value = 600
$ws.Range("K108").Value = 300
$ws.Range("L108").Value = 400
$ws.Range("M108").Value = 367
$ws.Range("P108").Value = 367
$ws.Range("D109").Value = 44510
$ws.Range("H109").Value = 'Paine'
$ws.Range("J109").Value = 2500
$ws.Range("K109").Value = 80
$ws.Range("L109").Value = 80
$ws.Range("M109").Value = 80
$ws.Range("P109").Value = 80
$ws.Range("H110").Value = 'Camote'
$ws.Range("J110").Value = 800
$ws.Range("K110").Value = 400
$ws.Range("L110").Value = 400
$ws.Range("M110").Value = 400
$ws.Range("P110").Value = 400
$ws.Range("D111").Value = 44463
$ws.Range("H111").Value = 'Paine'
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 130
$ws.Range("L111").Value = 130
$ws.Range("M111").Value = 130
$ws.Range("P111").Value = 130
$ws.Range("H112").Value = 'Camote'
$ws.Range("J112").Value = 900
$ws.Range("K112").Value = 400
$ws.Range("L112").Value = 400
$ws.Range("M112").Value = 400
$ws.Range("P112").Value = 400
$ws.Range("D113").Value = 44484
$ws.Range("H113").Value = 'Paine'
$ws.Range("I113").Value = '1a (guarda)'
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 100
$ws.Range("L113").Value = 100
$ws.Range("M113").Value = 100
$ws.Range("P113").Value = 100
$ws.Range("D114").Value = 44284
$ws.Range("I114").Value = '1a (cosecha)'
$ws.Range("J114").Value = 800
$ws.Range("D115").Value = 44376
$ws.Range("I115").Value = '1a (guarda)'
$ws.Range("J115").Value = 900
$ws.Range("K115").Value = 200
$ws.Range("L115").Value = 200
$ws.Range("M115").Value = 200
$ws.Range("P115").Value = 200
$ws.Range("D116").Value = 44242
$ws.Range("J116").Value = 800
$ws.Range("D117").Value = 44572
$ws.Range("I117").Value = '1a nueva(o)'
$ws.Range("J117").Value = 900
$ws.Range("K117").Value = 300
$ws.Range("L117").Value = 300
$ws.Range("M117").Value = 300
$ws.Range("O117").Value = 'Región del Maule'
$ws.Range("P117").Value = 300
$ws.Range("I118").Value = '1a (guarda)'
$ws.Range("K118").Value = 500
$ws.Range("L118").Value = 500
$ws.Range("M118").Value = 500
$ws.Range("O118").Value = 'Región de O''Higgins'
$ws.Range("P118").Value = 500
$ws.Range("H119").Value = 'Camote'
$ws.Range("I119").Value = '1a nueva(o)'
$ws.Range("J119").Value = 800
$ws.Range("K119").Value = 650
$ws.Range("L119").Value = 650
$ws.Range("M119").Value = 650
$ws.Range("O119").Value = 'Perú'
$ws.Range("P119").Value = 650
$ws.Range("D120").Value = 44516
$ws.Range("H120").Value = 'Paine'
$ws.Range("I120").Value = '1a (guarda)'
$ws.Range("J120").Value = 2500
$ws.Range("K120").Value = 80
$ws.Range("L120").Value = 80
$ws.Range("M120").Value = 80
$ws.Range("P120").Value = 80
$ws.Range("D121").Value = 44257
$ws.Range("I121").Value = '1a nueva(o)'
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 300
$ws.Range("L121").Value = 300
$ws.Range("M121").Value = 300
$ws.Range("P121").Value = 300
$ws.Range("D122").Value = 44279
$ws.Range("I122").Value = '1a (cosecha)'
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 200
$ws.Range("L122").Value = 200
$ws.Range("M122").Value = 200
$ws.Range("O122").Value = 'Región del Maule'
$ws.Range("P122").Value = 200
$ws.Range("D123").Value = 44552
$ws.Range("I123").Value = '1a nueva(o)'
$ws.Range("K123").Value = 650
$ws.Range("L123").Value = 650
$ws.Range("M123").Value = 650
$ws.Range("O123").Value = 'Región de O''Higgins'
$ws.Range("P123").Value = 650
$ws.Range("D124").Value = 44322
$ws.Range("I124").Value = '1a (guarda)'
$ws.Range("K124").Value = 280
$ws.Range("L124").Value = 280
$ws.Range("M124").Value = 280
$ws.Range("P124").Value = 280
$ws.Range("D125").Value = 44218
$ws.Range("I125").Value = '1a nueva(o)'
$ws.Range("D126").Value = 44400
$ws.Range("K126").Value = 350
$ws.Range("L126").Value = 350
$ws.Range("M126").Value = 350
$ws.Range("P126").Value = 350
$ws.Range("D127").Value = 44426
$ws.Range("K127").Value = 430
$ws.Range("L127").Value = 430
$ws.Range("M127").Value = 430
$ws.Range("P127").Value = 430
$ws.Range("D128").Value = 44431
$ws.Range("I128").Value = '1a (guarda)'
$ws.Range("K128").Value = 600
$ws.Range("L128").Value = 600
$ws.Range("M128").Value = 600
$ws.Range("P128").Value = 600
$ws.Range("D129").Value = 44211
$ws.Range("K129").Value = 300
$ws.Range("L129").Value = 300
$ws.Range("M129").Value = 300
$ws.Range("P129").Value = 300
$ws.Range("D130").Value = 44217
$ws.Range("K130").Value = 350
$ws.Range("L130").Value = 350
$ws.Range("M130").Value = 350
$ws.Range("O130").Value = 'Región del Maule'
$ws.Range("P130").Value = 350
$ws.Range("H131").Value = 'Camote'
$ws.Range("I131").Value = '1a nueva(o)'
$ws.Range("J131").Value = 800
$ws.Range("K131").Value = 600
$ws.Range("L131").Value = 600
$ws.Range("M131").Value = 600
$ws.Range("O131").Value = 'Perú'
$ws.Range("P131").Value = 600
$ws.Range("D132").Value = 44512
$ws.Range("H132").Value = 'Paine'
$ws.Range("I132").Value = '1a (guarda)'
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 80
$ws.Range("L132").Value = 80
$ws.Range("M132").Value = 80
$ws.Range("P132").Value = 80
$ws.Range("D133").Value = 44193
$ws.Range("J133").Value = 1500
$ws.Range("D134").Value = 44200
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 400
$ws.Range("L134").Value = 400
$ws.Range("M134").Value = 400
$ws.Range("P134").Value = 400
$ws.Range("D135").Value = 44203
$ws.Range("K135").Value = 350
$ws.Range("L135").Value = 350
$ws.Range("M135").Value = 350
$ws.Range("O135").Value = 'Región del Maule'
$ws.Range("P135").Value = 350
$ws.Range("D136").Value = 44545
$ws.Range("K136").Value = 600
$ws.Range("L136").Value = 600
$ws.Range("M136").Value = 600
$ws.Range("O136").Value = 'Región de O''Higgins'
$ws.Range("P136").Value = 600
$ws.Range("D137").Value = 44253
$ws.Range("H137").Value = 'Camote'
$ws.Range("I137").Value = '1a nueva(o)'
$ws.Range("J137").Value = 800
$ws.Range("K137").Value = 300
$ws.Range("L137").Value = 300
$ws.Range("M137").Value = 300
$ws.Range("P137").Value = 300
$ws.Range("D138").Value = 44504
$ws.Range("H138").Value = 'Paine'
$ws.Range("I138").Value = '1a (guarda)'
$ws.Range("J138").Value = 2500
$ws.Range("K138").Value = 80
$ws.Range("L138").Value = 80
$ws.Range("M138").Value = 80
$ws.Range("P138").Value = 80
$ws.Range("D139").Value = 44573
$ws.Range("I139").Value = '1a nueva(o)'
$ws.Range("J139").Value = 800
$ws.Range("K139").Value = 300
$ws.Range("L139").Value = 300
$ws.Range("M139").Value = 300
$ws.Range("P139").Value = 300
$ws.Range("D140").Value = 44280
$ws.Range("I140").Value = '1a (guarda)'
$ws.Range("J140").Value = 900
$ws.Range("K140").Value = 200
$ws.Range("L140").Value = 200
$ws.Range("M140").Value = 200
$ws.Range("P140").Value = 200
$ws.Range("D141").Value = 44187
$ws.Range("I141").Value = '1a nueva(o)'
$ws.Range("J141").Value = 1200
$ws.Range("K141").Value = 450
$ws.Range("L141").Value = 450
$ws.Range("M141").Value = 450
$ws.Range("P141").Value = 450
$ws.Range("D142").Value = 44391
$ws.Range("I142").Value = '1a (guarda)'
$ws.Range("J142").Value = 900
$ws.Range("K142").Value = 300
$ws.Range("L142").Value = 300
$ws.Range("M142").Value = 300
$ws.Range("P142").Value = 300
$ws.Range("D143").Value = 44202
$ws.Range("I143").Value = '1a nueva(o)'
$ws.Range("J143").Value = 800
$ws.Range("K143").Value = 350
$ws.Range("L143").Value = 350
$ws.Range("M143").Value = 350
$ws.Range("P143").Value = 350
$ws.Range("H144").Value = 'Camote'
$ws.Range("J144").Value = 900
$ws.Range("K144").Value = 280
$ws.Range("L144").Value = 280
$ws.Range("M144").Value = 280
$ws.Range("P144").Value = 280
$ws.Range("D145").Value = 44371
$ws.Range("H145").Value = 'Paine'
$ws.Range("I145").Value = '1a (guarda)'
$ws.Range("J145").Value = 1200
$ws.Range("K145").Value = 150
$ws.Range("L145").Value = 150
$ws.Range("M145").Value = 150
$ws.Range("P145").Value = 150
$ws.Range("D146").Value = 44249
$ws.Range("D147").Value = 44225
$ws.Range("H147").Value = 'Camote'
$ws.Range("I147").Value = '1a nueva(o)'
$ws.Range("J147").Value = 800
$ws.Range("K147").Value = 300
$ws.Range("L147").Value = 300
$ws.Range("M147").Value = 300
$ws.Range("P147").Value = 300
$ws.Range("D148").Value = 44445
$ws.Range("H148").Value = 'Paine'
$ws.Range("J148").Value = 2000
$ws.Range("K148").Value = 150
$ws.Range("L148").Value = 150
$ws.Range("M148").Value = 150
$ws.Range("P148").Value = 150
$ws.Range("D149").Value = 44330
$ws.Range("D150").Value = 44328
$ws.Range("I150").Value = '1a (guarda)'
$ws.Range("K150").Value = 280
$ws.Range("L150").Value = 280
$ws.Range("M150").Value = 280
$ws.Range("O150").Value = 'Región del Maule'
$ws.Range("P150").Value = 280
$ws.Range("H151").Value = 'Camote'
$ws.Range("I151").Value = '1a nueva(o)'
$ws.Range("J151").Value = 800
$ws.Range("K151").Value = 700
$ws.Range("L151").Value = 700
$ws.Range("M151").Value = 700
$ws.Range("O151").Value = 'Perú'
$ws.Range("P151").Value = 700
$ws.Range("D152").Value = 44526
$ws.Range("H152").Value = 'Paine'
$ws.Range("I152").Value = '1a (guarda)'
$ws.Range("J152").Value = 2000
$ws.Range("K152").Value = 80
$ws.Range("L152").Value = 80
$ws.Range("M152").Value = 80
$ws.Range("P152").Value = 80
$ws.Range("D153").Value = 44250
$ws.Range("H153").Value = 'Camote'
$ws.Range("I153").Value = '1a nueva(o)'
$ws.Range("J153").Value = 800
$ws.Range("K153").Value = 300
$ws.Range("L153").Value = 300
$ws.Range("M153").Value = 300
$ws.Range("P153").Value = 300
$ws.Range("D154").Value = 44285
$ws.Range("H154").Value = 'Sin especificar'
$ws.Range("J154").Value = 900
$ws.Range("K154").Value = 200
$ws.Range("L154").Value = 200
$ws.Range("M154").Value = 200
$ws.Range("P154").Value = 200
$ws.Range("H155").Value = 'Camote'
$ws.Range("J155").Value = 800
$ws.Range("K155").Value = 250
$ws.Range("L155").Value = 250
$ws.Range("M155").Value = 250
$ws.Range("P155").Value = 250
$ws.Range("D156").Value = 44264
$ws.Range("H156").Value = 'Paine'
$ws.Range("I156").Value = '1a (cosecha)'
$ws.Range("J156").Value = 1200
$ws.Range("K156").Value = 150
$ws.Range("L156").Value = 150
$ws.Range("M156").Value = 150
$ws.Range("P156").Value = 150
$ws.Range("D157").Value = 44221
$ws.Range("J157").Value = 1300
$ws.Range("K157").Value = 330
$ws.Range("L157").Value = 350
$ws.Range("M157").Value = 338
$ws.Range("O157").Value = 'Región del Maule'
$ws.Range("P157").Value = 338
$ws.Range("D158").Value = 44523
$ws.Range("I158").Value = '1a nueva(o)'
$ws.Range("K158").Value = 700
$ws.Range("L158").Value = 700
$ws.Range("M158").Value = 700
$ws.Range("O158").Value = 'Perú'
$ws.Range("P158").Value = 700
$ws.Range("D159").Value = 44399
$ws.Range("K159").Value = 350
$ws.Range("L159").Value = 350
$ws.Range("M159").Value = 350
$ws.Range("P159").Value = 350
$ws.Range("D160").Value = 44441
$ws.Range("J160").Value = 900
$ws.Range("K160").Value = 800
$ws.Range("L160").Value = 800
$ws.Range("M160").Value = 800
$ws.Range("P160").Value = 800
$ws.Range("H161").Value = 'Camote'
$ws.Range("J161").Value = 1000
$ws.Range("K161").Value = 200
$ws.Range("L161").Value = 200
$ws.Range("M161").Value = 200
$ws.Range("P161").Value = 200
$ws.Range("D162").Value = 44372
$ws.Range("H162").Value = 'Paine'
$ws.Range("I162").Value = '1a (guarda)'
$ws.Range("J162").Value = 1500
$ws.Range("K162").Value = 140
$ws.Range("L162").Value = 140
$ws.Range("M162").Value = 140
$ws.Range("P162").Value = 140
$ws.Range("D163").Value = 44186
$ws.Range("J163").Value = 800
$ws.Range("K163").Value = 450
$ws.Range("L163").Value = 450
$ws.Range("M163").Value = 450
$ws.Range("P163").Value = 450
$ws.Range("D164").Value = 44176
$ws.Range("I164").Value = '1a nueva(o)'
$ws.Range("J164").Value = 900
$ws.Range("K164").Value = 850
$ws.Range("L164").Value = 850
$ws.Range("M164").Value = 850
$ws.Range("P164").Value = 850
$ws.Range("H165").Value = 'Camote'
$ws.Range("J165").Value = 800
$ws.Range("K165").Value = 200
$ws.Range("L165").Value = 200
$ws.Range("M165").Value = 200
$ws.Range("P165").Value = 200
$ws.Range("D166").Value = 44278
$ws.Range("H166").Value = 'Paine'
$ws.Range("I166").Value = '1a (cosecha)'
$ws.Range("J166").Value = 1200
$ws.Range("K166").Value = 150
$ws.Range("L166").Value = 150
$ws.Range("M166").Value = 150
$ws.Range("P166").Value = 150
$ws.Range("H167").Value = 'Camote'
$ws.Range("J167").Value = 900
$ws.Range("K167").Value = 700
$ws.Range("L167").Value = 700
$ws.Range("M167").Value = 700
$ws.Range("P167").Value = 700
$ws.Range("D168").Value = 44438
$ws.Range("H168").Value = 'Paine'
$ws.Range("J168").Value = 2000
$ws.Range("K168").Value = 170
$ws.Range("L168").Value = 170
$ws.Range("M168").Value = 170
$ws.Range("P168").Value = 170
$ws.Range("D169").Value = 44300
$ws.Range("K169").Value = 200
$ws.Range("L169").Value = 200
$ws.Range("M169").Value = 200
$ws.Range("P169").Value = 200
$ws.Range("H170").Value = 'Camote'
$ws.Range("J170").Value = 900
$ws.Range("K170").Value = 600
$ws.Range("L170").Value = 600
$ws.Range("M170").Value = 600
$ws.Range("P170").Value = 600
$ws.Range("D171").Value = 44453
$ws.Range("H171").Value = 'Paine'
$ws.Range("J171").Value = 2000
$ws.Range("K171").Value = 150
$ws.Range("L171").Value = 150
$ws.Range("M171").Value = 150
$ws.Range("P171").Value = 150
$ws.Range("H172").Value = 'Camote'
$ws.Range("J172").Value = 1000
$ws.Range("K172").Value = 600
$ws.Range("L172").Value = 600
$ws.Range("M172").Value = 600
$ws.Range("P172").Value = 600
$ws.Range("D173").Value = 44449
$ws.Range("J173").Value = 1200
$ws.Range("K173").Value = 140
$ws.Range("L173").Value = 140
$ws.Range("M173").Value = 140
$ws.Range("P173").Value = 140
$ws.Range("D174").Value = 44468
$ws.Range("H174").Value = 'Paine'
$ws.Range("I174").Value = '1a (guarda)'
$ws.Range("J174").Value = 1500
$ws.Range("K174").Value = 150
$ws.Range("L174").Value = 150
$ws.Range("M174").Value = 150
$ws.Range("P174").Value = 150
$ws.Range("H175").Value = 'Camote'
$ws.Range("J175").Value = 900
$ws.Range("K175").Value = 250
$ws.Range("L175").Value = 250
$ws.Range("M175").Value = 250
$ws.Range("P175").Value = 250
$ws.Range("D176").Value = 44272
$ws.Range("H176").Value = 'Paine'
$ws.Range("I176").Value = '1a (cosecha)'
$ws.Range("J176").Value = 1200
$ws.Range("K176").Value = 150
$ws.Range("L176").Value = 150
$ws.Range("M176").Value = 150
$ws.Range("O176").Value = 'Región del Maule'
$ws.Range("P176").Value = 150
$ws.Range("D177").Value = 44529
$ws.Range("I177").Value = '1a nueva(o)'
$ws.Range("J177").Value = 800
$ws.Range("K177").Value = 700
$ws.Range("L177").Value = 700
$ws.Range("M177").Value = 700
$ws.Range("O177").Value = 'Perú'
$ws.Range("P177").Value = 700
$ws.Range("D178").Value = 44306
$ws.Range("D179").Value = 44299
$ws.Range("J179").Value = 900
$ws.Range("K179").Value = 200
$ws.Range("L179").Value = 200
$ws.Range("M179").Value = 200
$ws.Range("P179").Value = 200
$ws.Range("D180").Value = 44428
$ws.Range("J180").Value = 800
$ws.Range("K180").Value = 430
$ws.Range("L180").Value = 430
$ws.Range("M180").Value = 430
$ws.Range("P180").Value = 430
$ws.Range("D181").Value = 44321
$ws.Range("K181").Value = 280
$ws.Range("L181").Value = 280
$ws.Range("M181").Value = 280
$ws.Range("P181").Value = 280
$ws.Range("D182").Value = 44302
$ws.Range("I182").Value = '1a (guarda)'
$ws.Range("J182").Value = 900
$ws.Range("K182").Value = 200
$ws.Range("L182").Value = 200
$ws.Range("M182").Value = 200
$ws.Range("P182").Value = 200

# --- New row 183 (appended at the end, holding the former last record) ---
$ws.Range("A183").Value = 5
$ws.Range("B183").Value = 'Macroferia Regional de Talca'
$ws.Range("C183").Value = 'Maule'
$ws.Range("D183").Value = 44274
$ws.Range("E183").Value = 7
$ws.Range("F183").Value = 100112045
$ws.Range("G183").Value = 'Zapallo'
$ws.Range("H183").Value = 'Camote'
$ws.Range("I183").Value = '1a (cosecha)'
$ws.Range("J183").Value = 800
$ws.Range("K183").Value = 250
$ws.Range("L183").Value = 250
$ws.Range("M183").Value = 250
$ws.Range("N183").Value = '$/kilo (volumen en unidades)'
$ws.Range("O183").Value = 'Región del Maule'
$ws.Range("P183").Value = 250
$ws.Range("Q183").Value = 1
$ws.Range("R183").Value = 'Hortaliza'

# Match the date style (YYYY-MM-DD HH:MM:SS number format) used by the other Fecha cells
$ws.Range("D183").NumberFormat = $ws.Range("D182").NumberFormat
